$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in new "Asistido Copilot" hours for Bloque 1 rows (column D, rows 6-10)
$ws.Range("D6").Value = 2
$ws.Range("D7").Value = 2
$ws.Range("D8").Value = 4
$ws.Range("D9").Value = 4
$ws.Range("D10").Value = 2

# Update the view: scroll/selection moves to D10 (no more special top-left cell)
$ws.Range("D10").Select()
